# Configuracion IPv4 e IPv6
# Adds a new "transport" worksheet (TCP/UDP comparison + port ranges table)
# as the last tab in the workbook, mirroring the authored commit.

$wb = $excel.ActiveWorkbook

# Capture the current last sheet (IPv6) BEFORE adding the new sheet, since
# a freshly added sheet is inserted before the currently active sheet by
# default.
$lastExisting = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add($null, $lastExisting)
$ws.Name = "transport"

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "PROTOCOLO"
$ws.Range("B1").Value = "UDP"
$ws.Range("C1").Value = "TCP"
$ws.Range("E1").Value = "TIPO"
$ws.Range("F1").Value = "INICIO"
$ws.Range("G1").Value = "FIN"

# --- Significado ------------------------------------------------------------
$ws.Range("A2").Value = "SIGNIFICADO"
$ws.Range("B2").Value = "USER DATAGRAM PROTOCOL"
$ws.Range("C2").Value = "TRANSPORT CONTROL PROTOCOL"
$ws.Range("E2").Value = "CONOCIDO"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1023

# --- Caracteristica ---------------------------------------------------------
$ws.Range("A3").Value = "CARACTERISTICA"
$ws.Range("B3").Value = "MAXIMO ESFUERZO"
$ws.Range("C3").Value = "CONFIABLE (ACK)"
$ws.Range("E3").Value = "REGISTRADO"
$ws.Range("F3").Value = 1024
$ws.Range("G3").Value = 49151

# --- PDU ----------------------------------------------------------------
$ws.Range("A4").Value = "PDU"
$ws.Range("B4").Value = "DATAGRAMA"
$ws.Range("C4").Value = "SEGMENTO"
$ws.Range("E4").Value = "DINAMICO"
$ws.Range("F4").Value = 49152
$ws.Range("G4").Value = 65535

# --- Cabecera ----------------------------------------------------------
$ws.Range("A5").Value = "CABECERA"
$ws.Range("B5").Value = "8 BYTES"
$ws.Range("C5").Value = "20 BYTES"

# --- Aplicaciones --------------------------------------------------------
$ws.Range("A6").Value = "APLICACIONES"
$ws.Range("B6").Value = "DNS - DHCP - TFTP"
$ws.Range("C6").Value = "POP3 - IMAP - SMTP - SSH - HTTP - HTTPS - FTP - SFTP"

# --- Cosmetics: best-fit column widths, zoom, selection --------------------
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null

$ws.Range("C6").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130
